$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: insert the two new "Hub" columns right before the "Marketplace"
# column (currently column T). Doing this BEFORE the earlier insert (below)
# reproduces the shared-string insertion order seen in the target workbook
# (Hub ID / HUB Name were registered before Anniversary Date / Term Delta /
# Term Period).
# ---------------------------------------------------------------------------
$ws.Range("T1:U1").EntireColumn.Insert()
$ws.Range("T1").EntireColumn.ColumnWidth = 19.17
$ws.Range("T1").EntireColumn.OutlineLevel = 1
$ws.Range("U1").EntireColumn.ColumnWidth = 45
$ws.Range("U1").EntireColumn.OutlineLevel = 0
$ws.Range("T1").Value = "Hub ID"
$ws.Range("U1").Value = "HUB Name"

# ---------------------------------------------------------------------------
# Step 2: insert the three new columns (Anniversary Date, Term Delta,
# Term Period) right after "Subscription External ID" (column B), before the
# old "Customer ID" column.
# ---------------------------------------------------------------------------
$ws.Range("C1:E1").EntireColumn.Insert()
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 25
$ws.Range("C1:E1").EntireColumn.OutlineLevel = 1
$ws.Range("C1").Value = "Anniversary Date"
$ws.Range("D1").Value = "Term Delta"
$ws.Range("E1").Value = "Term Period"

# ---------------------------------------------------------------------------
# Step 3: update the worksheet AutoFilter range to match the new layout.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$discard = $ws.Range("K1:AE1").AutoFilter()

# ---------------------------------------------------------------------------
# Step 4: update the workbook-level _FilterDatabase defined name so it keeps
# pointing at the (now shifted) autofilter range.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Data!`$K`$1:`$AE`$1"
    }
}

# ---------------------------------------------------------------------------
# Step 5: update the active selection shown in the sheet view.
# ---------------------------------------------------------------------------
$discard2 = $ws.Range("C1:E1").Select()
